# Generate Report for Handoff
# Adds two new localized files (31a93d9a... and 4d8f17aa...) with status
# "Ready for handoff" to the Overview / zh-cn / de-de sheets & tables.

$wb = $excel.ActiveWorkbook

$repoBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1ec3096b9d04b8c70211aa28e7c3772e1701bd18/e2e/"

$file1 = "31a93d9a-6c2d-4f95-9781-35c4823c83b0.md"
$file2 = "4d8f17aa-978e-4fa6-a9bf-edb8ba655f6c.md"

$ho1zh = "31a93d9a-6c2d-4f95-9781-35c4823c83b0.fe0d7e0e606b72b60e0670cc6481891540c8ab1b.zh-cn.xlf"
$ho2zh = "4d8f17aa-978e-4fa6-a9bf-edb8ba655f6c.f3621c206cbc6d5cd2d74b72e86bdd88c4a4c3bd.zh-cn.xlf"
$ho1de = "31a93d9a-6c2d-4f95-9781-35c4823c83b0.fe0d7e0e606b72b60e0670cc6481891540c8ab1b.de-de.xlf"
$ho2de = "4d8f17aa-978e-4fa6-a9bf-edb8ba655f6c.f3621c206cbc6d5cd2d74b72e86bdd88c4a4c3bd.de-de.xlf"

$statusReady = "Ready for handoff"
$hoDateZh = "2016-10-20 08:10:17"
$hoDateDe = "2016-10-20 08:10:29"
$hbDateEmpty = "0001-01-01 00:00:00"

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A4").Value = $file1
$wsOverview.Range("B4").Value = "e2e\" + $file1
$wsOverview.Range("C4").Value = ".md"
$wsOverview.Range("D4").Value = ""
$wsOverview.Range("E4").Value = $statusReady
$wsOverview.Range("F4").Value = $statusReady
$wsOverview.Range("G4").Value = $hoDateDe

$wsOverview.Range("A5").Value = $file2
$wsOverview.Range("B5").Value = "e2e\" + $file2
$wsOverview.Range("C5").Value = ".md"
$wsOverview.Range("D5").Value = ""
$wsOverview.Range("E5").Value = $statusReady
$wsOverview.Range("F5").Value = $statusReady
$wsOverview.Range("G5").Value = $hoDateDe

$wsOverview.Hyperlinks.Add($wsOverview.Range("B4"), $repoBase + $file1, $null, $null, "e2e\" + $file1)
$wsOverview.Hyperlinks.Add($wsOverview.Range("B5"), $repoBase + $file2, $null, $null, "e2e\" + $file2)
$wsOverview.Range("B4").Style = "HyperLink"
$wsOverview.Range("B5").Style = "HyperLink"

$wsOverview.Range("G4:G5").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Columns.Item(5).ColumnWidth = 17.2159881591797
$wsOverview.Columns.Item(6).ColumnWidth = 17.2159881591797

$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G5"))

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A4").Value = $file1
$wsZh.Range("B4").Value = ".md"
$wsZh.Range("C4").Value = $statusReady
$wsZh.Range("D4").Value = "e2e"
$wsZh.Range("E4").Value = "ht"
$wsZh.Range("F4").Value = "'False"
$wsZh.Range("G4").Value = $ho1zh
$wsZh.Range("H4").Value = $hoDateZh
$wsZh.Range("I4").Value = ""
$wsZh.Range("J4").Value = ""
$wsZh.Range("K4").Value = $hbDateEmpty
$wsZh.Range("L4").Value = ""
$wsZh.Range("M4").Value = "'True"
$wsZh.Range("N4").Value = ""
$wsZh.Range("O4").Value = "'False"
$wsZh.Range("P4").Value = ""

$wsZh.Range("A5").Value = $file2
$wsZh.Range("B5").Value = ".md"
$wsZh.Range("C5").Value = $statusReady
$wsZh.Range("D5").Value = "e2e"
$wsZh.Range("E5").Value = "ht"
$wsZh.Range("F5").Value = "'False"
$wsZh.Range("G5").Value = $ho2zh
$wsZh.Range("H5").Value = $hoDateZh
$wsZh.Range("I5").Value = ""
$wsZh.Range("J5").Value = ""
$wsZh.Range("K5").Value = $hbDateEmpty
$wsZh.Range("L5").Value = ""
$wsZh.Range("M5").Value = "'True"
$wsZh.Range("N5").Value = ""
$wsZh.Range("O5").Value = "'False"
$wsZh.Range("P5").Value = ""

$wsZh.Hyperlinks.Add($wsZh.Range("A4"), $repoBase + $file1, $null, $null, $file1)
$wsZh.Hyperlinks.Add($wsZh.Range("A5"), $repoBase + $file2, $null, $null, $file2)
$wsZh.Range("A4").Style = "HyperLink"
$wsZh.Range("A5").Style = "HyperLink"

$wsZh.Range("H4:H5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("K4:K5").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsZh.Columns.Item(3).ColumnWidth = 17.2159881591797

$loZh = $wsZh.ListObjects.Item(1)
$loZh.Resize($wsZh.Range("A1:P5"))

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A4").Value = $file1
$wsDe.Range("B4").Value = ".md"
$wsDe.Range("C4").Value = $statusReady
$wsDe.Range("D4").Value = "e2e"
$wsDe.Range("E4").Value = "ht"
$wsDe.Range("F4").Value = "'False"
$wsDe.Range("G4").Value = $ho1de
$wsDe.Range("H4").Value = $hoDateDe
$wsDe.Range("I4").Value = ""
$wsDe.Range("J4").Value = ""
$wsDe.Range("K4").Value = $hbDateEmpty
$wsDe.Range("L4").Value = ""
$wsDe.Range("M4").Value = "'True"
$wsDe.Range("N4").Value = ""
$wsDe.Range("O4").Value = "'False"
$wsDe.Range("P4").Value = ""

$wsDe.Range("A5").Value = $file2
$wsDe.Range("B5").Value = ".md"
$wsDe.Range("C5").Value = $statusReady
$wsDe.Range("D5").Value = "e2e"
$wsDe.Range("E5").Value = "ht"
$wsDe.Range("F5").Value = "'False"
$wsDe.Range("G5").Value = $ho2de
$wsDe.Range("H5").Value = $hoDateDe
$wsDe.Range("I5").Value = ""
$wsDe.Range("J5").Value = ""
$wsDe.Range("K5").Value = $hbDateEmpty
$wsDe.Range("L5").Value = ""
$wsDe.Range("M5").Value = "'True"
$wsDe.Range("N5").Value = ""
$wsDe.Range("O5").Value = "'False"
$wsDe.Range("P5").Value = ""

$wsDe.Hyperlinks.Add($wsDe.Range("A4"), $repoBase + $file1, $null, $null, $file1)
$wsDe.Hyperlinks.Add($wsDe.Range("A5"), $repoBase + $file2, $null, $null, $file2)
$wsDe.Range("A4").Style = "HyperLink"
$wsDe.Range("A5").Style = "HyperLink"

$wsDe.Range("H4:H5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("K4:K5").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsDe.Columns.Item(3).ColumnWidth = 17.2159881591797

$loDe = $wsDe.ListObjects.Item(1)
$loDe.Resize($wsDe.Range("A1:P5"))

Write-Host "Done applying handoff report rows."
